# Sprint 2 backlog update
# Reassign the "Play Minesweeper" user story (row 5) from Ali Cooper to
# Caleb Ljunggren, and leave the selection where the editor's cursor
# ended up (H5) after making the change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5").Value = "Caleb Ljunggren"

$ws.Range("H5").Select() | Out-Null
